# "Se procesan de nuevo los datos con las nuevas dimensiones curadas"
#
# The metadata sheet documents, per data column (row 1 = column name),
# how that column should be interpreted:
#   row 2 -> iaest-dimension:* / iaest-measure:* / sdmx-dimension:* annotation
#   row 3 -> annotation kind ("dim" or "medida")
#   row 4 -> datatype (skos:Concept / xsd:int / URI-<something>)
#   row 5 -> mapping workbook file (only present for skos:Concept dimensions)
#
# With the newly curated dimensions:
#   - "tipo-de-vivienda-principal" (column F) stops being a curated
#     dimension and becomes a plain numeric measure.
#   - "tipo-de-vivienda" (column J) stops being a curated dimension and
#     becomes a plain numeric measure.
#   - "municipio-nombre" (column K) stops being a measure and becomes a
#     geographic reference-area dimension, just like provincia-nombre
#     (column L) and comarca-nombre (column M), with its own URI column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F: tipo-de-vivienda-principal -> now a measure (was a dimension)
$ws.Range("F2").Value = "iaest-measure:tipo-de-vivienda-principal"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"

# Column J: tipo-de-vivienda -> now a measure (was a dimension)
$ws.Range("J2").Value = "iaest-measure:tipo-de-vivienda"
$ws.Range("J3").Value = "medida"
$ws.Range("J4").Value = "xsd:int"

# Column K: municipio-nombre -> now a refArea dimension (was a measure),
# matching provincia-nombre / comarca-nombre, with its own URI column.
$ws.Range("K2").Value = "sdmx-dimension:refArea"
$ws.Range("K3").Value = "dim"
$ws.Range("K4").Value = "URI-Municipio"

# F and J no longer reference a mapping workbook (they are measures now),
# so their row-5 mapping-file cells are removed entirely.
$ws.Cells.Item(5, 6).Clear()
$ws.Cells.Item(5, 10).Clear()
